$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 14; existing rows 14-104 shift down to 15-105.
$ws.Rows("14:14").Insert()

# Populate the newly inserted row 14 with its data.
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C14").Value = "Ñuble"
$ws.Range("D14").Value = 44950
$ws.Range("E14").Value = 16
$ws.Range("F14").Value = 100112022
$ws.Range("G14").Value = "Arveja Verde"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 21000
$ws.Range("L14").Value = 21000
$ws.Range("M14").Value = 21000
$ws.Range("N14").Value = "$/saco 25 kilos"
$ws.Range("O14").Value = "Región de Ñuble"
$ws.Range("P14").Value = 840
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"

# Keep the date style/number format used by the rest of column D.
$ws.Range("D14").NumberFormat = $ws.Range("D15").NumberFormat
